$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 1230.1
$ws.Range("I62").Value = 1270.7142
$ws.Range("J62").Value = 1135.3334
$ws.Range("K62").Value = 1270.7142
$ws.Range("L62").Value = 1135.3334
$ws.Range("M62").Value = -646.7141999999999
$ws.Range("N62").Value = -2383.3334
$ws.Range("H65").Value = 1230.1
$ws.Range("I65").Value = 1270.7142
$ws.Range("J65").Value = 1135.3334
$ws.Range("K65").Value = 6353.571
$ws.Range("L65").Value = 5676.666999999999
$ws.Range("M65").Value = -3233.571
$ws.Range("N65").Value = -11916.667
$ws.Range("H80").Value = 695.65
$ws.Range("I80").Value = 597.1539
$ws.Range("K80").Value = 1791.4617
$ws.Range("M80").Value = -793.4617000000001
$ws.Range("H83").Value = 695.65
$ws.Range("I83").Value = 597.1539
$ws.Range("K83").Value = 5374.3851
$ws.Range("M83").Value = -382.3851000000004
$ws.Range("H86").Value = 3538.35
$ws.Range("I86").Value = 4277.923
$ws.Range("K86").Value = 4277.923
$ws.Range("M86").Value = -3154.923
$ws.Range("H89").Value = 3538.35
$ws.Range("I89").Value = 4277.923
$ws.Range("K89").Value = 21389.615
$ws.Range("M89").Value = -15773.615
$ws.Range("H107").Value = 25006420
$ws.Range("J107").Value = 17933.334
$ws.Range("L107").Value = 17933.334
$ws.Range("N107").Value = -21773.334
$ws.Range("H113").Value = 3230.875
$ws.Range("I113").Value = 2027.8572
$ws.Range("J113").Value = 4166.5557
$ws.Range("K113").Value = 2027.8572
$ws.Range("L113").Value = 4166.5557
$ws.Range("M113").Value = 1226.1428
$ws.Range("N113").Value = -10674.5557
$ws.Range("H132").Value = 2122.6667
$ws.Range("I132").Value = 2122.6667
$ws.Range("K132").Value = 6368.000100000001
$ws.Range("M132").Value = -3838.000100000001
$ws.Range("H136").Value = 28280
$ws.Range("J136").Value = 28280
$ws.Range("L136").Value = 28280
$ws.Range("N136").Value = -38480
$ws.Range("H137").Value = 1766.7778
$ws.Range("I137").Value = 1714.4286
$ws.Range("K137").Value = 5143.2858
$ws.Range("M137").Value = -2593.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6332.915
$ws.Range("I32").Value = 4874.457
$ws.Range("J32").Value = 10586.75
$ws.Range("K32").Value = 4874.457
$ws.Range("L32").Value = 10586.75
$ws.Range("M32").Value = -4587.457
$ws.Range("N32").Value = -11160.75
$ws.Range("H132").Value = 3831.8845
$ws.Range("I132").Value = 1361.125
$ws.Range("J132").Value = 7785.1
$ws.Range("K132").Value = 4083.375
$ws.Range("L132").Value = 23355.3
$ws.Range("M132").Value = -1553.375
$ws.Range("N132").Value = -28415.3
$ws.Range("H140").Value = 35161.54
$ws.Range("J140").Value = 35161.54
$ws.Range("L140").Value = 35161.54
$ws.Range("N140").Value = -45521.54

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4834.8613
$ws.Range("I134").Value = 5652.115
$ws.Range("J134").Value = 2710
$ws.Range("K134").Value = 16956.345
$ws.Range("L134").Value = 8130
$ws.Range("M134").Value = -14421.345
$ws.Range("N134").Value = -13200
$ws.Range("H137").Value = 38569.668
$ws.Range("J137").Value = 38000
$ws.Range("L137").Value = 38000
$ws.Range("N137").Value = -48200

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 5129584
$ws.Range("I16").Value = 6994296
$ws.Range("J16").Value = 1625
$ws.Range("K16").Value = 6994296
$ws.Range("L16").Value = 1625
$ws.Range("M16").Value = -6994009
$ws.Range("N16").Value = -2199
$ws.Range("H31").Value = 6813
$ws.Range("I31").Value = 2055.1177
$ws.Range("J31").Value = 12205.267
$ws.Range("K31").Value = 2055.1177
$ws.Range("L31").Value = 12205.267
$ws.Range("M31").Value = -1760.1177
$ws.Range("N31").Value = -12795.267
$ws.Range("H34").Value = 6813
$ws.Range("I34").Value = 2055.1177
$ws.Range("J34").Value = 12205.267
$ws.Range("K34").Value = 2055.1177
$ws.Range("L34").Value = 12205.267
$ws.Range("M34").Value = -1853.1177
$ws.Range("N34").Value = -12609.267
$ws.Range("H58").Value = 1552.7693
$ws.Range("I58").Value = 1165
$ws.Range("J58").Value = 1795.125
$ws.Range("K58").Value = 1165
$ws.Range("L58").Value = 1795.125
$ws.Range("M58").Value = -962
$ws.Range("N58").Value = -2201.125
$ws.Range("H99").Value = 8931811
$ws.Range("I99").Value = 1027.875
$ws.Range("K99").Value = 1027.875
$ws.Range("M99").Value = 470.125
$ws.Range("H113").Value = 5129584
$ws.Range("I113").Value = 6994296
$ws.Range("J113").Value = 1625
$ws.Range("K113").Value = 6994296
$ws.Range("L113").Value = 1625
$ws.Range("M113").Value = -6992126
$ws.Range("N113").Value = -5965
$ws.Range("H126").Value = 8931811
$ws.Range("I126").Value = 1027.875
$ws.Range("K126").Value = 3083.625
$ws.Range("M126").Value = -613.625
$ws.Range("H134").Value = 2649.2942
$ws.Range("I134").Value = 3018.5
$ws.Range("J134").Value = 1763.2
$ws.Range("K134").Value = 9055.5
$ws.Range("L134").Value = 5289.6
$ws.Range("M134").Value = -6520.5
$ws.Range("N134").Value = -10359.6
$ws.Range("H136").Value = 1552.7693
$ws.Range("I136").Value = 1165
$ws.Range("J136").Value = 1795.125
$ws.Range("K136").Value = 3495
$ws.Range("L136").Value = 5385.375
$ws.Range("M136").Value = -945
$ws.Range("N136").Value = -10485.375

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 311.33334
$ws.Range("I99").Value = 311.33334
$ws.Range("K99").Value = 934.0000200000001
$ws.Range("M99").Value = 1311.99998
$ws.Range("H131").Value = 1370816.8
$ws.Range("I131").Value = 5883095.5
$ws.Range("J131").Value = 1017.8036
$ws.Range("K131").Value = 17649286.5
$ws.Range("L131").Value = 3053.4108
$ws.Range("M131").Value = -17644246.5
$ws.Range("N131").Value = -13133.4108
$ws.Range("H138").Value = 30307.5
$ws.Range("I138").Value = 30307.5
$ws.Range("K138").Value = 90922.5
$ws.Range("M138").Value = -85782.5
$ws.Range("H139").Value = 5895.484
$ws.Range("I139").Value = 34686.668
$ws.Range("J139").Value = 2810.7144
$ws.Range("K139").Value = 104060.004
$ws.Range("L139").Value = 8432.143199999999
$ws.Range("M139").Value = -98920.00399999999
$ws.Range("N139").Value = -18712.1432
$ws.Range("H140").Value = 2734.4375
$ws.Range("I140").Value = 2613.7273
$ws.Range("K140").Value = 7841.1819
$ws.Range("M140").Value = -2661.1819
$ws.Range("H141").Value = 17609.182
$ws.Range("J141").Value = 9557.143
$ws.Range("L141").Value = 28671.429
$ws.Range("N141").Value = -39031.429

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 5469241
$ws.Range("J122").Value = 6252014
$ws.Range("L122").Value = 18756042
$ws.Range("N122").Value = -18760942

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 76924990
$ws.Range("I40").Value = 111112890
$ws.Range("J40").Value = 2213.75
$ws.Range("K40").Value = 111112890
$ws.Range("L40").Value = 2213.75
$ws.Range("M40").Value = -111112754
$ws.Range("N40").Value = -2485.75
$ws.Range("H68").Value = 71430620
$ws.Range("I68").Value = 1669.6666
$ws.Range("K68").Value = 1669.6666
$ws.Range("M68").Value = -920.6666
$ws.Range("H71").Value = 71430620
$ws.Range("I71").Value = 1669.6666
$ws.Range("K71").Value = 8348.333000000001
$ws.Range("M71").Value = -4604.333000000001
$ws.Range("H112").Value = 38220.8
$ws.Range("J112").Value = 38220.8
$ws.Range("L112").Value = 38220.8
$ws.Range("N112").Value = -41174.8
$ws.Range("H122").Value = 3879722.5
$ws.Range("I122").Value = 4764215
$ws.Range("J122").Value = 1668491.6
$ws.Range("K122").Value = 14292645
$ws.Range("L122").Value = 5005474.800000001
$ws.Range("M122").Value = -14290195
$ws.Range("N122").Value = -5010374.800000001
$ws.Range("H132").Value = 12750990
$ws.Range("I132").Value = 18062006
$ws.Range("J132").Value = 4554.3
$ws.Range("K132").Value = 54186018
$ws.Range("L132").Value = 13662.9
$ws.Range("M132").Value = -54183488
$ws.Range("N132").Value = -18722.9
$ws.Range("H135").Value = 46714.5
$ws.Range("J135").Value = 46714.5
$ws.Range("L135").Value = 46714.5
$ws.Range("N135").Value = -56854.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2733.5
$ws.Range("I96").Value = 2050
$ws.Range("J96").Value = 3075.25
$ws.Range("K96").Value = 2050
$ws.Range("L96").Value = 3075.25
$ws.Range("M96").Value = -677
$ws.Range("N96").Value = -5821.25
$ws.Range("H122").Value = 2112.7778
$ws.Range("I122").Value = 2019.1666
$ws.Range("K122").Value = 6057.4998
$ws.Range("M122").Value = -3607.4998
